$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 22153.555  # H40
$ws.Cells.Item(40, 10).Value = 27876.8  # J40
$ws.Cells.Item(40, 12).Value = 27876.8  # L40
$ws.Cells.Item(40, 14).Value = -28226.8  # N40

$ws.Cells.Item(86, 8).Value = 200003410  # H86
$ws.Cells.Item(86, 9).Value = 250003140  # I86
$ws.Cells.Item(86, 11).Value = 250003140  # K86
$ws.Cells.Item(86, 13).Value = -250002017  # M86

$ws.Cells.Item(89, 8).Value = 200003410  # H89
$ws.Cells.Item(89, 9).Value = 250003140  # I89
$ws.Cells.Item(89, 11).Value = 1250015700  # K89
$ws.Cells.Item(89, 13).Value = -1250010084  # M89

$ws.Cells.Item(103, 8).Value = 1471.25  # H103
$ws.Cells.Item(103, 9).Value = 1394.2  # I103
$ws.Cells.Item(103, 11).Value = 4182.6  # K103
$ws.Cells.Item(103, 13).Value = -3596.6  # M103

$ws.Cells.Item(115, 8).Value = 142864690  # H115
$ws.Cells.Item(115, 10).Value = 0  # J115
$ws.Cells.Item(115, 12).Value = 0  # L115
$ws.Cells.Item(115, 14).ClearContents()  # N115

$ws.Cells.Item(132, 8).Value = 4862.0835  # H132
$ws.Cells.Item(132, 9).Value = 3844.5854  # I132
$ws.Cells.Item(132, 11).Value = 11533.7562  # K132
$ws.Cells.Item(132, 13).Value = -9003.7562  # M132

$ws.Cells.Item(137, 8).Value = 2529.5635  # H137
$ws.Cells.Item(137, 9).Value = 2491.6216  # I137
$ws.Cells.Item(137, 10).Value = 2570.853  # J137
$ws.Cells.Item(137, 11).Value = 7474.864799999999  # K137
$ws.Cells.Item(137, 12).Value = 7712.559  # L137
$ws.Cells.Item(137, 13).Value = -4924.864799999999  # M137
$ws.Cells.Item(137, 14).Value = -12812.559  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 9122.885  # H132
$ws.Cells.Item(132, 9).Value = 8103  # I132
$ws.Cells.Item(132, 10).Value = 11891.143  # J132
$ws.Cells.Item(132, 11).Value = 24309  # K132
$ws.Cells.Item(132, 12).Value = 35673.429  # L132
$ws.Cells.Item(132, 13).Value = -21779  # M132
$ws.Cells.Item(132, 14).Value = -40733.429  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3840.0588  # H134
$ws.Cells.Item(134, 9).Value = 2979.0222  # I134
$ws.Cells.Item(134, 11).Value = 8937.0666  # K134
$ws.Cells.Item(134, 13).Value = -6402.0666  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3047.158  # H16
$ws.Cells.Item(16, 9).Value = 2411.4167  # I16
$ws.Cells.Item(16, 11).Value = 2411.4167  # K16
$ws.Cells.Item(16, 13).Value = -2124.4167  # M16

$ws.Cells.Item(31, 8).Value = 2089.7273  # H31
$ws.Cells.Item(31, 9).Value = 1915.2  # I31
$ws.Cells.Item(31, 10).Value = 2358.2307  # J31
$ws.Cells.Item(31, 11).Value = 1915.2  # K31
$ws.Cells.Item(31, 12).Value = 2358.2307  # L31
$ws.Cells.Item(31, 13).Value = -1620.2  # M31
$ws.Cells.Item(31, 14).Value = -2948.2307  # N31

$ws.Cells.Item(34, 8).Value = 2089.7273  # H34
$ws.Cells.Item(34, 9).Value = 1915.2  # I34
$ws.Cells.Item(34, 10).Value = 2358.2307  # J34
$ws.Cells.Item(34, 11).Value = 1915.2  # K34
$ws.Cells.Item(34, 12).Value = 2358.2307  # L34
$ws.Cells.Item(34, 13).Value = -1713.2  # M34
$ws.Cells.Item(34, 14).Value = -2762.2307  # N34

$ws.Cells.Item(58, 8).Value = 5417.375  # H58
$ws.Cells.Item(58, 9).Value = 5980.0527  # I58
$ws.Cells.Item(58, 11).Value = 5980.0527  # K58
$ws.Cells.Item(58, 13).Value = -5777.0527  # M58

$ws.Cells.Item(86, 8).Value = 3710359.2  # H86
$ws.Cells.Item(86, 9).Value = 6674208  # I86
$ws.Cells.Item(86, 11).Value = 6674208  # K86
$ws.Cells.Item(86, 13).Value = -6673085  # M86

$ws.Cells.Item(89, 8).Value = 3710359.2  # H89
$ws.Cells.Item(89, 9).Value = 6674208  # I89
$ws.Cells.Item(89, 11).Value = 33371040  # K89
$ws.Cells.Item(89, 13).Value = -33365424  # M89

$ws.Cells.Item(99, 8).Value = 10968.0625  # H99
$ws.Cells.Item(99, 9).Value = 4206.3335  # I99
$ws.Cells.Item(99, 10).Value = 15025.1  # J99
$ws.Cells.Item(99, 11).Value = 4206.3335  # K99
$ws.Cells.Item(99, 12).Value = 15025.1  # L99
$ws.Cells.Item(99, 13).Value = -2708.3335  # M99
$ws.Cells.Item(99, 14).Value = -18021.1  # N99

$ws.Cells.Item(105, 8).Value = 1828.7646  # H105
$ws.Cells.Item(105, 9).Value = 1828.7646  # I105
$ws.Cells.Item(105, 10).Value = 0  # J105
$ws.Cells.Item(105, 11).Value = 1828.7646  # K105
$ws.Cells.Item(105, 12).Value = 0  # L105
$ws.Cells.Item(105, 13).Value = -81.76459999999997  # M105
$ws.Cells.Item(105, 14).ClearContents()  # N105

$ws.Cells.Item(107, 8).Value = 694.6061  # H107
$ws.Cells.Item(107, 9).Value = 679.3226  # I107
$ws.Cells.Item(107, 10).Value = 931.5  # J107
$ws.Cells.Item(107, 11).Value = 679.3226  # K107
$ws.Cells.Item(107, 12).Value = 931.5  # L107
$ws.Cells.Item(107, 13).Value = 1240.6774  # M107
$ws.Cells.Item(107, 14).Value = -4771.5  # N107

$ws.Cells.Item(113, 8).Value = 3047.158  # H113
$ws.Cells.Item(113, 9).Value = 2411.4167  # I113
$ws.Cells.Item(113, 11).Value = 2411.4167  # K113
$ws.Cells.Item(113, 13).Value = -241.4167000000002  # M113

$ws.Cells.Item(122, 8).Value = 5546.1665  # H122
$ws.Cells.Item(122, 9).Value = 1642.7142  # I122
$ws.Cells.Item(122, 11).Value = 4928.142599999999  # K122
$ws.Cells.Item(122, 13).Value = -2478.142599999999  # M122

$ws.Cells.Item(126, 8).Value = 10968.0625  # H126
$ws.Cells.Item(126, 9).Value = 4206.3335  # I126
$ws.Cells.Item(126, 10).Value = 15025.1  # J126
$ws.Cells.Item(126, 11).Value = 12619.0005  # K126
$ws.Cells.Item(126, 12).Value = 45075.3  # L126
$ws.Cells.Item(126, 13).Value = -10149.0005  # M126
$ws.Cells.Item(126, 14).Value = -50015.3  # N126

$ws.Cells.Item(132, 8).Value = 6692.3335  # H132
$ws.Cells.Item(132, 9).Value = 7337.8096  # I132
$ws.Cells.Item(132, 10).Value = 4433.1665  # J132
$ws.Cells.Item(132, 11).Value = 22013.4288  # K132
$ws.Cells.Item(132, 12).Value = 13299.4995  # L132
$ws.Cells.Item(132, 13).Value = -19483.4288  # M132
$ws.Cells.Item(132, 14).Value = -18359.4995  # N132

$ws.Cells.Item(134, 8).Value = 5574.0244  # H134
$ws.Cells.Item(134, 9).Value = 4801.3125  # I134
$ws.Cells.Item(134, 10).Value = 8321.444  # J134
$ws.Cells.Item(134, 11).Value = 14403.9375  # K134
$ws.Cells.Item(134, 12).Value = 24964.332  # L134
$ws.Cells.Item(134, 13).Value = -11868.9375  # M134
$ws.Cells.Item(134, 14).Value = -30034.332  # N134

$ws.Cells.Item(136, 8).Value = 5417.375  # H136
$ws.Cells.Item(136, 9).Value = 5980.0527  # I136
$ws.Cells.Item(136, 11).Value = 17940.1581  # K136
$ws.Cells.Item(136, 13).Value = -15390.1581  # M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(29, 8).Value = 106  # H29
$ws.Cells.Item(29, 10).Value = 86  # J29
$ws.Cells.Item(29, 12).Value = 258  # L29
$ws.Cells.Item(29, 14).Value = -812  # N29

$ws.Cells.Item(63, 8).Value = 8513.615  # H63
$ws.Cells.Item(63, 9).Value = 8513.615  # I63
$ws.Cells.Item(63, 10).Value = 0  # J63
$ws.Cells.Item(63, 11).Value = 25540.845  # K63
$ws.Cells.Item(63, 12).Value = 0  # L63
$ws.Cells.Item(63, 13).Value = -24791.845  # M63
$ws.Cells.Item(63, 14).ClearContents()  # N63

$ws.Cells.Item(66, 8).Value = 8513.615  # H66
$ws.Cells.Item(66, 9).Value = 8513.615  # I66
$ws.Cells.Item(66, 10).Value = 0  # J66
$ws.Cells.Item(66, 11).Value = 76622.535  # K66
$ws.Cells.Item(66, 12).Value = 0  # L66
$ws.Cells.Item(66, 13).Value = -72878.535  # M66
$ws.Cells.Item(66, 14).ClearContents()  # N66

$ws.Cells.Item(68, 8).Value = 1822.5143  # H68
$ws.Cells.Item(68, 9).Value = 1499  # I68
$ws.Cells.Item(68, 11).Value = 4497  # K68
$ws.Cells.Item(68, 13).Value = -3686  # M68

$ws.Cells.Item(71, 8).Value = 1822.5143  # H71
$ws.Cells.Item(71, 9).Value = 1499  # I71
$ws.Cells.Item(71, 11).Value = 13491  # K71
$ws.Cells.Item(71, 13).Value = -9435  # M71

$ws.Cells.Item(80, 8).Value = 3996.3333  # H80
$ws.Cells.Item(80, 9).Value = 3995  # I80
$ws.Cells.Item(80, 11).Value = 11985  # K80
$ws.Cells.Item(80, 13).Value = -11049  # M80

$ws.Cells.Item(83, 8).Value = 3996.3333  # H83
$ws.Cells.Item(83, 9).Value = 3995  # I83
$ws.Cells.Item(83, 11).Value = 35955  # K83
$ws.Cells.Item(83, 13).Value = -31275  # M83

$ws.Cells.Item(113, 8).Value = 1749.8667  # H113
$ws.Cells.Item(113, 9).Value = 3379.6667  # I113
$ws.Cells.Item(113, 10).Value = 1342.4166  # J113
$ws.Cells.Item(113, 11).Value = 10139.0001  # K113
$ws.Cells.Item(113, 12).Value = 4027.2498  # L113
$ws.Cells.Item(113, 13).Value = -7969.000100000001  # M113
$ws.Cells.Item(113, 14).Value = -8367.2498  # N113

$ws.Cells.Item(122, 8).Value = 871.1539  # H122
$ws.Cells.Item(122, 9).Value = 630.5333000000001  # I122
$ws.Cells.Item(122, 10).Value = 1199.2727  # J122
$ws.Cells.Item(122, 11).Value = 5674.7997  # K122
$ws.Cells.Item(122, 12).Value = 10793.4543  # L122
$ws.Cells.Item(122, 13).Value = -3224.7997  # M122
$ws.Cells.Item(122, 14).Value = -15693.4543  # N122

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 3883.9  # H70
$ws.Cells.Item(70, 9).Value = 2558.25  # I70
$ws.Cells.Item(70, 10).Value = 5872.375  # J70
$ws.Cells.Item(70, 11).Value = 2558.25  # K70
$ws.Cells.Item(70, 12).Value = 5872.375  # L70
$ws.Cells.Item(70, 13).Value = -2288.25  # M70
$ws.Cells.Item(70, 14).Value = -6412.375  # N70

$ws.Cells.Item(73, 8).Value = 3883.9  # H73
$ws.Cells.Item(73, 9).Value = 2558.25  # I73
$ws.Cells.Item(73, 10).Value = 5872.375  # J73
$ws.Cells.Item(73, 11).Value = 2558.25  # K73
$ws.Cells.Item(73, 12).Value = 5872.375  # L73
$ws.Cells.Item(73, 13).Value = -1622.25  # M73
$ws.Cells.Item(73, 14).Value = -7744.375  # N73

$ws.Cells.Item(122, 8).Value = 2035.0857  # H122
$ws.Cells.Item(122, 9).Value = 1900.381  # I122
$ws.Cells.Item(122, 11).Value = 5701.143  # K122
$ws.Cells.Item(122, 13).Value = -3251.143  # M122

$ws.Cells.Item(132, 8).Value = 3915.3447  # H132
$ws.Cells.Item(132, 9).Value = 3967.8147  # I132
$ws.Cells.Item(132, 10).Value = 3207  # J132
$ws.Cells.Item(132, 11).Value = 11903.4441  # K132
$ws.Cells.Item(132, 12).Value = 9621  # L132
$ws.Cells.Item(132, 13).Value = -9373.444100000001  # M132
$ws.Cells.Item(132, 14).Value = -14681  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 27136.244  # H132
$ws.Cells.Item(132, 9).Value = 31314.086  # I132
$ws.Cells.Item(132, 10).Value = 2765.5  # J132
$ws.Cells.Item(132, 11).Value = 93942.258  # K132
$ws.Cells.Item(132, 12).Value = 8296.5  # L132
$ws.Cells.Item(132, 13).Value = -91412.258  # M132
$ws.Cells.Item(132, 14).Value = -13356.5  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 8361.069  # H122
$ws.Cells.Item(122, 9).Value = 4315.2104  # I122
$ws.Cells.Item(122, 11).Value = 12945.6312  # K122
$ws.Cells.Item(122, 13).Value = -10495.6312  # M122

$ws.Cells.Item(132, 8).Value = 2374.2932  # H132
$ws.Cells.Item(132, 9).Value = 2323.302  # I132
$ws.Cells.Item(132, 11).Value = 6969.906000000001  # K132
$ws.Cells.Item(132, 13).Value = -4439.906000000001  # M132

$ws.Cells.Item(138, 8).Value = 98994.8  # H138
$ws.Cells.Item(138, 10).Value = 98994.8  # J138
$ws.Cells.Item(138, 12).Value = 98994.8  # L138
$ws.Cells.Item(138, 14).Value = -109274.8  # N138
